$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.591.88'
$ws.Range('E2').Value = '  -3.39%  '

$ws.Range('D3').Value = '2.091.25'
$ws.Range('E3').Value = '  -1.11%  '

$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '342.25'
$ws.Range('E5').Value = '  -2.17%  '

$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  -0.25%  '

$ws.Range('D7').Value = '0.5123'
$ws.Range('E7').Value = '  -2.44%  '

$ws.Range('D8').Value = '0.4382'
$ws.Range('E8').Value = '  -2.87%  '

$ws.Range('D9').Value = '52.68'
$ws.Range('E9').Value = '  -2.12%  '

$ws.Range('D10').Value = '0.09060'
$ws.Range('E10').Value = '  +0.39%  '

$ws.Range('D11').Value = '1.169'
$ws.Range('E11').Value = '  -0.28%  '

$ws.Range('D12').Value = '24.79'
$ws.Range('E12').Value = '  +1.21%  '

$ws.Range('D13').Value = '2.086.40'
$ws.Range('E13').Value = '  -1.85%  '

$ws.Range('D14').Value = '6.758'
$ws.Range('E14').Value = '  -1.05%  '

$ws.Range('D15').Value = '8.199'
$ws.Range('E15').Value = '  +2.10%  '

$ws.Range('D16').Value = '99.97'
$ws.Range('E16').Value = '  -1.98%  '

$ws.Range('E17').Value = '  -0.35%  '

$ws.Range('D18').Value = '0.00001143'
$ws.Range('E18').Value = '  -2.19%  '

$ws.Range('D19').Value = '20.98'
$ws.Range('E19').Value = '  +8.10%  '

$ws.Range('D20').Value = '0.06631'
$ws.Range('E20').Value = '  -1.17%  '

$ws.Range('D21').Value = '1.008'
$ws.Range('E21').Value = '  -0.13%  '

$ws.Range('D22').Value = '6.162'
$ws.Range('E22').Value = '  -2.26%  '

$ws.Range('D23').Value = '29.616.30'
$ws.Range('E23').Value = '  -3.55%  '

$ws.Range('D24').Value = '12.59'
$ws.Range('E24').Value = '  -1.92%  '

$ws.Range('D25').Value = '2.299'
$ws.Range('E25').Value = '  -3.75%  '

$ws.Range('D26').Value = '2.342.26'
$ws.Range('E26').Value = '  -1.33%  '

$ws.Range('D27').Value = '21.72'
$ws.Range('E27').Value = '  -2.92%  '

$ws.Range('D28').Value = '163.04'
$ws.Range('E28').Value = '  -1.35%  '

$ws.Range('D29').Value = '2.512'
$ws.Range('E29').Value = '  -1.07%  '

$ws.Range('E30').Value = '  -3.09%  '

$ws.Range('D31').Value = '1.126'
$ws.Range('E31').Value = '  -5.43%  '

$ws.Range('D32').Value = '0.1043'
$ws.Range('E32').Value = '  -3.17%  '

$ws.Range('D33').Value = '1.625'
$ws.Range('E33').Value = '  -1.14%  '

$ws.Range('D34').Value = '6.122'
$ws.Range('E34').Value = '  -4.06%  '

$ws.Range('D35').Value = '3.955'
$ws.Range('E35').Value = '  -1.66%  '

$ws.Range('D36').Value = '6.018'
$ws.Range('E36').Value = '  +1.58%  '

$ws.Range('D37').Value = '10.21'
$ws.Range('E37').Value = '  -1.13%  '

$ws.Range('D38').Value = '0.02561'
$ws.Range('E38').Value = '  -3.60%  '

$ws.Range('D39').Value = '0.06659'
$ws.Range('E39').Value = '  -2.69%  '

$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '12.36'
$ws.Range('E40').Value = '  -1.47%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.2227'
$ws.Range('E41').Value = '  -3.97%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6826'
$ws.Range('E42').Value = '  -0.92%  '

$ws.Range('D43').Value = '1.279'
$ws.Range('E43').Value = '  +0.52%  '

$ws.Range('D44').Value = '0.6632'
$ws.Range('E44').Value = '  +2.87%  '

$ws.Range('D45').Value = '14.05'
$ws.Range('E45').Value = '  -4.46%  '

$ws.Range('D46').Value = '2.285'
$ws.Range('E46').Value = '  -1.87%  '

$ws.Range('D47').Value = '3.604'
$ws.Range('E47').Value = '  -3.92%  '

$ws.Range('D48').Value = '1.215'
$ws.Range('E48').Value = '  -2.90%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '81.53'
$ws.Range('E49').Value = '  -1.13%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.00000000334'
$ws.Range('E50').Value = '  -6.21%  '

$ws.Range('D51').Value = '1.168'
$ws.Range('E51').Value = '  -2.00%  '
